$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying source data (a rolling sample of price records) was
# refreshed for the weekly update: almost every existing row (2-26) was
# reassigned to a different position using the same set of records (one of
# them, formerly row 18, also had its Variedad corrected from
# "Sin especificar" to "Magnum"), a new record was inserted as row 20, and
# the data set grew by one additional row, appended as row 27.

# Row 2
$ws.Cells.Item(2, 4).Value = 44279
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 100112031
$ws.Cells.Item(2, 7).Value = 'Poroto verde'
$ws.Cells.Item(2, 8).Value = 'Magnum'
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 100
$ws.Cells.Item(2, 11).Value = 28000
$ws.Cells.Item(2, 12).Value = 30000
$ws.Cells.Item(2, 13).Value = 29000
$ws.Cells.Item(2, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(2, 15).Value = 'Región del Maule'
$ws.Cells.Item(2, 16).Value = 1160
$ws.Cells.Item(2, 17).Value = 25
$ws.Cells.Item(2, 18).Value = 'Hortaliza'

# Row 3
$ws.Cells.Item(3, 4).Value = 44237
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 100112031
$ws.Cells.Item(3, 7).Value = 'Poroto verde'
$ws.Cells.Item(3, 8).Value = 'Sin especificar'
$ws.Cells.Item(3, 9).Value = 'Primera'
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 20000
$ws.Cells.Item(3, 12).Value = 22000
$ws.Cells.Item(3, 13).Value = 21000
$ws.Cells.Item(3, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(3, 15).Value = 'Región del Maule'
$ws.Cells.Item(3, 16).Value = 840
$ws.Cells.Item(3, 17).Value = 25
$ws.Cells.Item(3, 18).Value = 'Hortaliza'

# Row 4
$ws.Cells.Item(4, 4).Value = 44441
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 100112031
$ws.Cells.Item(4, 7).Value = 'Poroto verde'
$ws.Cells.Item(4, 8).Value = 'Magnum'
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 28000
$ws.Cells.Item(4, 12).Value = 29000
$ws.Cells.Item(4, 13).Value = 28500
$ws.Cells.Item(4, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(4, 15).Value = 'Perú'
$ws.Cells.Item(4, 16).Value = 1140
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = 'Hortaliza'

# Row 5
$ws.Cells.Item(5, 4).Value = 44253
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = 100112031
$ws.Cells.Item(5, 7).Value = 'Poroto verde'
$ws.Cells.Item(5, 8).Value = 'Magnum'
$ws.Cells.Item(5, 9).Value = 'Primera'
$ws.Cells.Item(5, 10).Value = 200
$ws.Cells.Item(5, 11).Value = 25000
$ws.Cells.Item(5, 12).Value = 26000
$ws.Cells.Item(5, 13).Value = 25500
$ws.Cells.Item(5, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5, 15).Value = 'Región del Maule'
$ws.Cells.Item(5, 16).Value = 1020
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = 'Hortaliza'

# Row 6
$ws.Cells.Item(6, 4).Value = 44323
$ws.Cells.Item(6, 5).Value = 8
$ws.Cells.Item(6, 6).Value = 100112031
$ws.Cells.Item(6, 7).Value = 'Poroto verde'
$ws.Cells.Item(6, 8).Value = 'Magnum'
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 100
$ws.Cells.Item(6, 11).Value = 20000
$ws.Cells.Item(6, 12).Value = 22000
$ws.Cells.Item(6, 13).Value = 21000
$ws.Cells.Item(6, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(6, 15).Value = 'Perú'
$ws.Cells.Item(6, 16).Value = 840
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = 'Hortaliza'

# Row 7
$ws.Cells.Item(7, 4).Value = 44342
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 100112031
$ws.Cells.Item(7, 7).Value = 'Poroto verde'
$ws.Cells.Item(7, 8).Value = 'Magnum'
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 28000
$ws.Cells.Item(7, 12).Value = 30000
$ws.Cells.Item(7, 13).Value = 29000
$ws.Cells.Item(7, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(7, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(7, 16).Value = 1160
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = 'Hortaliza'

# Row 8
$ws.Cells.Item(8, 4).Value = 44160
$ws.Cells.Item(8, 5).Value = 8
$ws.Cells.Item(8, 6).Value = 100112031
$ws.Cells.Item(8, 7).Value = 'Poroto verde'
$ws.Cells.Item(8, 8).Value = 'Magnum'
$ws.Cells.Item(8, 9).Value = 'Primera'
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 28000
$ws.Cells.Item(8, 12).Value = 30000
$ws.Cells.Item(8, 13).Value = 29000
$ws.Cells.Item(8, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(8, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 16).Value = 1160
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = 'Hortaliza'

# Row 9
$ws.Cells.Item(9, 4).Value = 44167
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 100112031
$ws.Cells.Item(9, 7).Value = 'Poroto verde'
$ws.Cells.Item(9, 8).Value = 'Sin especificar'
$ws.Cells.Item(9, 9).Value = 'Primera'
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 18000
$ws.Cells.Item(9, 12).Value = 19000
$ws.Cells.Item(9, 13).Value = 18500
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9, 15).Value = 'Región del Maule'
$ws.Cells.Item(9, 16).Value = 740
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = 'Hortaliza'

# Row 10
$ws.Cells.Item(10, 4).Value = 44272
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 100112031
$ws.Cells.Item(10, 7).Value = 'Poroto verde'
$ws.Cells.Item(10, 8).Value = 'Magnum'
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 100
$ws.Cells.Item(10, 11).Value = 22000
$ws.Cells.Item(10, 12).Value = 24000
$ws.Cells.Item(10, 13).Value = 23000
$ws.Cells.Item(10, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(10, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(10, 16).Value = 920
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = 'Hortaliza'

# Row 11
$ws.Cells.Item(11, 4).Value = 44399
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 100112031
$ws.Cells.Item(11, 7).Value = 'Poroto verde'
$ws.Cells.Item(11, 8).Value = 'Magnum'
$ws.Cells.Item(11, 9).Value = 'Primera'
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 20000
$ws.Cells.Item(11, 12).Value = 22000
$ws.Cells.Item(11, 13).Value = 21000
$ws.Cells.Item(11, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(11, 15).Value = 'Perú'
$ws.Cells.Item(11, 16).Value = 840
$ws.Cells.Item(11, 17).Value = 25
$ws.Cells.Item(11, 18).Value = 'Hortaliza'

# Row 12
$ws.Cells.Item(12, 4).Value = 44188
$ws.Cells.Item(12, 5).Value = 8
$ws.Cells.Item(12, 6).Value = 100112031
$ws.Cells.Item(12, 7).Value = 'Poroto verde'
$ws.Cells.Item(12, 8).Value = 'Magnum'
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 38000
$ws.Cells.Item(12, 12).Value = 40000
$ws.Cells.Item(12, 13).Value = 39000
$ws.Cells.Item(12, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(12, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(12, 16).Value = 1560
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = 'Hortaliza'

# Row 13
$ws.Cells.Item(13, 4).Value = 44230
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = 100112031
$ws.Cells.Item(13, 7).Value = 'Poroto verde'
$ws.Cells.Item(13, 8).Value = 'Magnum'
$ws.Cells.Item(13, 9).Value = 'Primera'
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 22000
$ws.Cells.Item(13, 12).Value = 24000
$ws.Cells.Item(13, 13).Value = 23000
$ws.Cells.Item(13, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(13, 15).Value = 'Región del Maule'
$ws.Cells.Item(13, 16).Value = 920
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = 'Hortaliza'

# Row 14
$ws.Cells.Item(14, 4).Value = 44244
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 100112031
$ws.Cells.Item(14, 7).Value = 'Poroto verde'
$ws.Cells.Item(14, 8).Value = 'Magnum'
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 16000
$ws.Cells.Item(14, 12).Value = 18000
$ws.Cells.Item(14, 13).Value = 17000
$ws.Cells.Item(14, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(14, 15).Value = 'Región del Maule'
$ws.Cells.Item(14, 16).Value = 680
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = 'Hortaliza'

# Row 15
$ws.Cells.Item(15, 4).Value = 44265
$ws.Cells.Item(15, 5).Value = 8
$ws.Cells.Item(15, 6).Value = 100112031
$ws.Cells.Item(15, 7).Value = 'Poroto verde'
$ws.Cells.Item(15, 8).Value = 'Magnum'
$ws.Cells.Item(15, 9).Value = 'Primera'
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 20000
$ws.Cells.Item(15, 12).Value = 22000
$ws.Cells.Item(15, 13).Value = 21000
$ws.Cells.Item(15, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(15, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 16).Value = 840
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = 'Hortaliza'

# Row 16
$ws.Cells.Item(16, 4).Value = 44447
$ws.Cells.Item(16, 5).Value = 8
$ws.Cells.Item(16, 6).Value = 100112031
$ws.Cells.Item(16, 7).Value = 'Poroto verde'
$ws.Cells.Item(16, 8).Value = 'Magnum'
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 100
$ws.Cells.Item(16, 11).Value = 37000
$ws.Cells.Item(16, 12).Value = 38000
$ws.Cells.Item(16, 13).Value = 37500
$ws.Cells.Item(16, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(16, 15).Value = 'Perú'
$ws.Cells.Item(16, 16).Value = 1500
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = 'Hortaliza'

# Row 17
$ws.Cells.Item(17, 4).Value = 44294
$ws.Cells.Item(17, 5).Value = 8
$ws.Cells.Item(17, 6).Value = 100112031
$ws.Cells.Item(17, 7).Value = 'Poroto verde'
$ws.Cells.Item(17, 8).Value = 'Magnum'
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 100
$ws.Cells.Item(17, 11).Value = 24000
$ws.Cells.Item(17, 12).Value = 25000
$ws.Cells.Item(17, 13).Value = 24500
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(17, 15).Value = 'Región del Maule'
$ws.Cells.Item(17, 16).Value = 980
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = 'Hortaliza'

# Row 18
$ws.Cells.Item(18, 4).Value = 44203
$ws.Cells.Item(18, 5).Value = 8
$ws.Cells.Item(18, 6).Value = 100112031
$ws.Cells.Item(18, 7).Value = 'Poroto verde'
$ws.Cells.Item(18, 8).Value = 'Magnum'
$ws.Cells.Item(18, 9).Value = 'Primera'
$ws.Cells.Item(18, 10).Value = 100
$ws.Cells.Item(18, 11).Value = 20000
$ws.Cells.Item(18, 12).Value = 22000
$ws.Cells.Item(18, 13).Value = 21000
$ws.Cells.Item(18, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(18, 15).Value = 'Región del Maule'
$ws.Cells.Item(18, 16).Value = 840
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = 'Hortaliza'

# Row 19
$ws.Cells.Item(19, 4).Value = 44335
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(19, 6).Value = 100112031
$ws.Cells.Item(19, 7).Value = 'Poroto verde'
$ws.Cells.Item(19, 8).Value = 'Magnum'
$ws.Cells.Item(19, 9).Value = 'Primera'
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 35000
$ws.Cells.Item(19, 12).Value = 36000
$ws.Cells.Item(19, 13).Value = 35500
$ws.Cells.Item(19, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(19, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(19, 16).Value = 1420
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = 'Hortaliza'

# Row 20
$ws.Cells.Item(20, 4).Value = 44475
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 100112031
$ws.Cells.Item(20, 7).Value = 'Poroto verde'
$ws.Cells.Item(20, 8).Value = 'Magnum'
$ws.Cells.Item(20, 9).Value = 'Primera'
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 44000
$ws.Cells.Item(20, 12).Value = 45000
$ws.Cells.Item(20, 13).Value = 44500
$ws.Cells.Item(20, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(20, 15).Value = 'Perú'
$ws.Cells.Item(20, 16).Value = 1780
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = 'Hortaliza'

# Row 21
$ws.Cells.Item(21, 4).Value = 44435
$ws.Cells.Item(21, 5).Value = 8
$ws.Cells.Item(21, 6).Value = 100112031
$ws.Cells.Item(21, 7).Value = 'Poroto verde'
$ws.Cells.Item(21, 8).Value = 'Magnum'
$ws.Cells.Item(21, 9).Value = 'Primera'
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 25000
$ws.Cells.Item(21, 12).Value = 26000
$ws.Cells.Item(21, 13).Value = 25500
$ws.Cells.Item(21, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(21, 15).Value = 'Perú'
$ws.Cells.Item(21, 16).Value = 1020
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = 'Hortaliza'

# Row 22
$ws.Cells.Item(22, 4).Value = 44433
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = 100112031
$ws.Cells.Item(22, 7).Value = 'Poroto verde'
$ws.Cells.Item(22, 8).Value = 'Magnum'
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 25000
$ws.Cells.Item(22, 12).Value = 26000
$ws.Cells.Item(22, 13).Value = 25500
$ws.Cells.Item(22, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(22, 15).Value = 'Perú'
$ws.Cells.Item(22, 16).Value = 1020
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = 'Hortaliza'

# Row 23
$ws.Cells.Item(23, 4).Value = 44468
$ws.Cells.Item(23, 5).Value = 8
$ws.Cells.Item(23, 6).Value = 100112031
$ws.Cells.Item(23, 7).Value = 'Poroto verde'
$ws.Cells.Item(23, 8).Value = 'Magnum'
$ws.Cells.Item(23, 9).Value = 'Primera'
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 31000
$ws.Cells.Item(23, 12).Value = 32000
$ws.Cells.Item(23, 13).Value = 31500
$ws.Cells.Item(23, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(23, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(23, 16).Value = 1260
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = 'Hortaliza'

# Row 24
$ws.Cells.Item(24, 4).Value = 44384
$ws.Cells.Item(24, 5).Value = 8
$ws.Cells.Item(24, 6).Value = 100112031
$ws.Cells.Item(24, 7).Value = 'Poroto verde'
$ws.Cells.Item(24, 8).Value = 'Sin especificar'
$ws.Cells.Item(24, 9).Value = 'Primera'
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 25000
$ws.Cells.Item(24, 12).Value = 26000
$ws.Cells.Item(24, 13).Value = 25500
$ws.Cells.Item(24, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(24, 15).Value = 'Perú'
$ws.Cells.Item(24, 16).Value = 1020
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = 'Hortaliza'

# Row 25
$ws.Cells.Item(25, 4).Value = 44363
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(25, 6).Value = 100112031
$ws.Cells.Item(25, 7).Value = 'Poroto verde'
$ws.Cells.Item(25, 8).Value = 'Magnum'
$ws.Cells.Item(25, 9).Value = 'Primera'
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 25000
$ws.Cells.Item(25, 12).Value = 26000
$ws.Cells.Item(25, 13).Value = 25500
$ws.Cells.Item(25, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(25, 15).Value = 'Perú'
$ws.Cells.Item(25, 16).Value = 1020
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = 'Hortaliza'

# Row 26
$ws.Cells.Item(26, 4).Value = 44321
$ws.Cells.Item(26, 5).Value = 8
$ws.Cells.Item(26, 6).Value = 100112031
$ws.Cells.Item(26, 7).Value = 'Poroto verde'
$ws.Cells.Item(26, 8).Value = 'Magnum'
$ws.Cells.Item(26, 9).Value = 'Primera'
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 24000
$ws.Cells.Item(26, 12).Value = 25000
$ws.Cells.Item(26, 13).Value = 24500
$ws.Cells.Item(26, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(26, 15).Value = 'Región del Maule'
$ws.Cells.Item(26, 16).Value = 980
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = 'Hortaliza'

# Row 27
$ws.Cells.Item(27, 4).Value = 44461
$ws.Cells.Item(27, 5).Value = 8
$ws.Cells.Item(27, 6).Value = 100112031
$ws.Cells.Item(27, 7).Value = 'Poroto verde'
$ws.Cells.Item(27, 8).Value = 'Sin especificar'
$ws.Cells.Item(27, 9).Value = 'Primera'
$ws.Cells.Item(27, 10).Value = 100
$ws.Cells.Item(27, 11).Value = 33000
$ws.Cells.Item(27, 12).Value = 34000
$ws.Cells.Item(27, 13).Value = 33500
$ws.Cells.Item(27, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(27, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(27, 16).Value = 1340
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = 'Hortaliza'

# Row 27 is a brand-new data row; fill in the Mercado ID / Mercado / Region
# columns (A-C), which are constant across every row in this sheet.
$ws.Cells.Item(27, 1).Value = 11
$ws.Cells.Item(27, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(27, 3).Value = 'Bíobío'

# Apply the same date number format used for the Fecha (D) column elsewhere
# in the sheet so the new row's date cell matches the existing formatting.
$ws.Cells.Item(27, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
